$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is updated: "Placeholder sprite characters" row becomes "Dude Monster" row,
# and the Note column changes from "Probably not used in the final version" to "Placeholder"
$ws.Range("A4").Value = "Dude Monster"
$ws.Range("B4").Value = "Craftpix.net"
$ws.Range("C4").Value = "https://free-game-assets.itch.io/free-tiny-hero-sprites-pixel-art"
$ws.Range("D4").Value = "Placeholder"

# Row 5 stays the PaperZD Tutorial credit (unchanged content)
$ws.Range("A5").Value = "PaperZD Tutorial"
$ws.Range("B5").Value = "LeafBranchGames"
$ws.Range("C5").Value = "https://www.youtube.com/watch?v=aWkgOr5U-zI"

# New row 6: Valla sprite credit
$ws.Range("A6").Value = "Valla"
$ws.Range("B6").Value = "Badim"
$ws.Range("C6").Value = "https://badim.itch.io/pixelart-valla"
$ws.Range("D6").Value = "Placeholder"

# New row 7: Slimes sprite credit
$ws.Range("A7").Value = "Slimes"
$ws.Range("B7").Value = "Shimmy"
$ws.Range("C7").Value = "https://shimyx.itch.io/6-direction-slimes-animation"
$ws.Range("D7").Value = "Placeholder"

# New row 8: 2D Topdown Tutorial credit
$ws.Range("A8").Value = "2D Topdown Tutorial"
$ws.Range("B8").Value = "Cobra Code"
$ws.Range("C8").Value = "https://www.youtube.com/watch?v=z1RMDMKcROQ"

# Update the selected cell to reflect where the author left off editing
$ws.Range("C9").Select()
